$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.390.37'
$ws.Range('E2').Value = '  -2.55%  '

$ws.Range('D3').Value = '3.075.78'
$ws.Range('E3').Value = '  -3.53%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '525.84'
$ws.Range('E5').Value = '  -2.06%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.75'
$ws.Range('E6').Value = '  -4.68%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.22%  '

$ws.Range('D8').Value = '3.074.84'
$ws.Range('E8').Value = '  -3.45%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.471'
$ws.Range('E9').Value = '  +4.41%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  +0.01%  '

$ws.Range('E11').Value = '  -4.66%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.405'
$ws.Range('E12').Value = '  +0.97%  '

$ws.Range('E13').Value = '  +2.07%  '

$ws.Range('D14').Value = '3.606.79'
$ws.Range('E14').Value = '  -3.55%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '25.20'
$ws.Range('E15').Value = '  -3.14%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000161'
$ws.Range('E16').Value = '  -5.22%  '

$ws.Range('D17').Value = '57.409.68'
$ws.Range('E17').Value = '  -2.53%  '

$ws.Range('D18').Value = '3.075.64'
$ws.Range('E18').Value = '  -3.66%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.89'
$ws.Range('E19').Value = '  -4.89%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.47'
$ws.Range('E20').Value = '  -3.98%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.84'
$ws.Range('E21').Value = '  -3.40%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '350.19'
$ws.Range('E22').Value = '  -2.58%  '

$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.06'
$ws.Range('E24').Value = '  -1.20%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.499'
$ws.Range('E25').Value = '  -3.80%  '

$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.165'
$ws.Range('E26').Value = '  -3.25%  '

$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.05%  '

$ws.Range('D28').Value = '0.0₃0851'
$ws.Range('E28').Value = '  -11.43%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.09'
$ws.Range('E30').Value = '  -6.47%  '

$ws.Range('E31').Value = '  -3.48%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.93'
$ws.Range('E32').Value = '  -9.83%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '21.04'
$ws.Range('E33').Value = '  -2.19%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.82'
$ws.Range('E34').Value = '  -2.44%  '

$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '158.53'
$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.13'
$ws.Range('E36').Value = '  -7.42%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.99'
$ws.Range('E37').Value = '  -5.32%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '25.41'
$ws.Range('E38').Value = '  -4.93%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.23'
$ws.Range('E39').Value = '  -6.74%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0660'
$ws.Range('E40').Value = '  -2.93%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.58'
$ws.Range('E41').Value = '  -5.99%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.01'
$ws.Range('E42').Value = '  -1.62%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.689'
$ws.Range('E43').Value = '  -3.31%  '

$ws.Range('D44').Value = '2.406.93'
$ws.Range('E44').Value = '  +1.28%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '36.79'
$ws.Range('E45').Value = '  -0.73%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('D47').Value = '3.115.89'
$ws.Range('E47').Value = '  -3.48%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0260'
$ws.Range('E48').Value = '  -4.78%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.00'
$ws.Range('E49').Value = '  -2.00%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.948'
$ws.Range('E50').Value = '  -7.83%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.35'
$ws.Range('E51').Value = '  -7.17%  '
